$d = $word.ActiveDocument

$pairs = @(
    @("431×8=", "549×4="),
    @("905×8=", "852×7="),
    @("528×2=", "793×5="),
    @("795×8=", "883×2="),
    @("384×8=", "231×6="),
    @("847×2=", "185×3="),
    @("718×2=", "594×7="),
    @("599×9=", "191×8="),
    @("515×8=", "809×8="),
    @("477×2=", "426×4="),
    @("670×4=", "778×3="),
    @("970×5=", "381×9="),
    @("868×5=", "875×4="),
    @("738×2=", "174×7="),
    @("357×6=", "113×6="),
    @("922×9=", "576×7="),
    @("650×6=", "656×3="),
    @("362×6=", "849×7="),
    @("757×8=", "679×6="),
    @("475×8=", "261×8="),
    @("818×3=", "543×2="),
    @("611×3=", "876×5="),
    @("559×2=", "640×6="),
    @("800×9=", "730×3="),
    @("485×8=", "329×5=")
)

foreach ($pair in $pairs) {
    $old = $pair[0]
    $new = $pair[1]
    $range = $d.Content
    $range.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}
